$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) are stored as text in the source data.
# Force text number-format first so Excel does not silently coerce numeric-
# looking strings (e.g. "92.50") into numbers and drop the trailing zero.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "41.645.54"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.468.04"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "318.44"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "92.50"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").Value = "0.552"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").Value = "0.0872"
$ws.Range("E10").Value = "  +10.04%  "
$ws.Range("D11").Value = "33.01"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "2.847.65"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "6.89"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "15.50"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "2.471.08"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("D18").Value = "41.582.79"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "70.85"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "11.30"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "240.71"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "24.86"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "36.52"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("D31").Value = "157.84"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "5.48"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +4.47%  "
$ws.Range("D38").Value = "2.90"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +8.13%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "1.989.26"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("D44").Value = "18.94"
$ws.Range("E44").Value = "  +1.79%  "
$ws.Range("D45").Value = "0.0283"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("D48").Value = "2.705.15"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "97.64"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").Value = "75.59"
$ws.Range("E50").Value = "  +5.85%  "
$ws.Range("D51").Value = "66.82"
$ws.Range("E51").Value = "  -0.24%  "
